$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text interpretation for the Price/Volume columns so numeric-looking
# strings (e.g. "92.64", "1.00") are preserved verbatim as text instead of
# being coerced into floating point numbers (which would drop trailing
# zeros / thousands separators).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "42.387.86"
$ws.Range("E2").Value = "  -8.45%  "
$ws.Range("D3").Value = "2.508.60"
$ws.Range("E3").Value = "  -3.68%  "
$ws.Range("D4").Value = "0.997"
$ws.Range("E4").Value = "  -0.33%  "
$ws.Range("D5").Value = "294.88"
$ws.Range("E5").Value = "  -3.74%  "
$ws.Range("D6").Value = "92.64"
$ws.Range("E6").Value = "  -6.82%  "
$ws.Range("D7").Value = "0.567"
$ws.Range("E7").Value = "  -5.48%  "
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").Value = "0.543"
$ws.Range("E9").Value = "  -5.79%  "
$ws.Range("D10").Value = "35.91"
$ws.Range("E10").Value = "  -8.56%  "
$ws.Range("D11").Value = "0.0795"
$ws.Range("E11").Value = "  -5.25%  "
$ws.Range("D12").Value = "7.60"
$ws.Range("E12").Value = "  -6.04%  "
$ws.Range("E13").Value = "  +0.24%  "
$ws.Range("D14").Value = "2.873.65"
$ws.Range("E14").Value = "  -4.26%  "
$ws.Range("D15").Value = "2.490.08"
$ws.Range("E15").Value = "  -4.63%  "
$ws.Range("D16").Value = "0.860"
$ws.Range("E16").Value = "  -6.08%  "
$ws.Range("D17").Value = "13.96"
$ws.Range("E17").Value = "  -6.22%  "
$ws.Range("D18").Value = "42.290.98"
$ws.Range("E18").Value = "  -8.79%  "
$ws.Range("D19").Value = "0.0₃0955"
$ws.Range("E19").Value = "  -4.97%  "
$ws.Range("D20").Value = "6.46"
$ws.Range("E20").Value = "  -3.16%  "
$ws.Range("D21").Value = "12.23"
$ws.Range("E21").Value = "  -5.18%  "
$ws.Range("D22").Value = "72.36"
$ws.Range("E22").Value = "  +1.62%  "
$ws.Range("D23").Value = "256.54"
$ws.Range("E23").Value = "  -5.63%  "
$ws.Range("D24").Value = "2.86"
$ws.Range("E24").Value = "  -5.16%  "
$ws.Range("D25").Value = "2.09"
$ws.Range("E25").Value = "  -3.07%  "
$ws.Range("D26").Value = "28.70"
$ws.Range("E26").Value = "  -2.18%  "
$ws.Range("E27").Value = "  +0.11%  "
$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").Value = "9.85"
$ws.Range("E28").Value = "  -6.57%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "2.20"
$ws.Range("E29").Value = "  -0.66%  "
$ws.Range("D30").Value = "36.53"
$ws.Range("E30").Value = "  -4.23%  "
$ws.Range("D31").Value = "5.97"
$ws.Range("E31").Value = "  -5.21%  "
$ws.Range("D32").Value = "3.44"
$ws.Range("E32").Value = "  -5.32%  "
$ws.Range("D33").Value = "150.68"
$ws.Range("E33").Value = "  -0.44%  "
$ws.Range("B34").Value = "EnergySwap"
$ws.Range("C34").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D34").Value = "27.57"
$ws.Range("E34").Value = "  +19.37%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "2.16"
$ws.Range("E35").Value = "  -2.85%  "
$ws.Range("B36").Value = "WEMIXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").Value = "2.69"
$ws.Range("E36").Value = "  -5.78%  "
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").Value = "0.0791"
$ws.Range("E37").Value = "  -4.91%  "
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").Value = "0.113"
$ws.Range("E38").Value = "  -7.32%  "
$ws.Range("D39").Value = "0.118"
$ws.Range("E39").Value = "  -3.90%  "
$ws.Range("D40").Value = "16.49"
$ws.Range("E40").Value = "  +4.36%  "
$ws.Range("D41").Value = "3.40"
$ws.Range("E41").Value = "  -5.01%  "
$ws.Range("D42").Value = "0.0306"
$ws.Range("E42").Value = "  -6.73%  "
$ws.Range("D43").Value = "3.81"
$ws.Range("E43").Value = "  -5.82%  "
$ws.Range("D44").Value = "2.003.89"
$ws.Range("E44").Value = "  -5.30%  "
$ws.Range("D45").Value = "0.995"
$ws.Range("E45").Value = "  -0.30%  "
$ws.Range("D46").Value = "85.15"
$ws.Range("E46").Value = "  -8.91%  "
$ws.Range("E47").Value = "  +2.44%  "
$ws.Range("D48").Value = "8.78"
$ws.Range("E48").Value = "  -7.69%  "
$ws.Range("D49").Value = "2.735.87"
$ws.Range("E49").Value = "  -4.25%  "
$ws.Range("D50").Value = "102.24"
$ws.Range("E50").Value = "  -5.54%  "
$ws.Range("D51").Value = "1.63"
$ws.Range("E51").Value = "  -7.83%  "

# Restore the default (unstyled) cell formatting so the edited cells keep
# the same style as before the edit.
$ws.Range("D2:E51").Style = "Normal"

